$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns H (x) and I (y)
$ws.Range("H1").Value = "x"
$ws.Range("I1").Value = "y"

# Map coordinate data per territory row
$xyData = @{
    2 = @(5, 35)
    3 = @(11, 32)
    4 = @(19, 31)
    5 = @(6, 26)
    6 = @(12, 28)
    7 = @(16, 27)
    8 = @(11, 24)
    9 = @(8, 22)
    10 = @(16, 17)
    11 = @(15, 12)
    12 = @(19, 14)
    13 = @(24, 14)
    14 = @(21, 12)
    15 = @(22, 9)
    16 = @(18, 6)
    17 = @(34, 32)
    18 = @(34, 28)
    19 = @(37, 30)
    20 = @(37, 27)
    21 = @(41, 30)
    22 = @(42, 26)
    23 = @(43, 35)
    24 = @(44, 31)
    25 = @(53, 34)
    26 = @(63, 34)
    27 = @(49, 30)
    28 = @(57, 30)
    29 = @(49, 26)
    30 = @(54, 27)
    31 = @(60, 27)
    32 = @(59, 24)
    33 = @(63, 23)
    34 = @(64.5, 26)
    35 = @(67.5, 26)
    36 = @(53, 23)
    37 = @(52, 20)
    38 = @(56, 22)
    39 = @(54, 19)
    40 = @(60, 19)
    41 = @(63, 15)
    42 = @(65, 18)
    43 = @(67, 8)
    44 = @(72, 5)
    45 = @(47, 24)
    46 = @(45, 22)
    47 = @(37, 22)
    48 = @(32, 19)
    49 = @(37, 18)
    50 = @(44, 16)
    51 = @(39, 13)
    52 = @(41, 8)
    53 = @(3, 22)
    54 = @(5, 15)
    55 = @(25, 26)
    56 = @(22, 18)
    57 = @(28, 8)
    58 = @(34, 13)
    59 = @(39, 25)
    60 = @(53, 10)
    61 = @(71, 27)
    62 = @(70, 19)
    63 = @(65, 3)
}

foreach ($r in $xyData.Keys) {
    $vals = $xyData[$r]
    $ws.Cells.Item([int]$r, 8).Value = $vals[0]
    $ws.Cells.Item([int]$r, 9).Value = $vals[1]
}

# Update view: scroll position and active selection
$ws.Range("I22").Select()
$excel.ActiveWindow.ScrollRow = 12
